# "Stage Class & Cover Sheet"
# Adds a Duration value to the existing Development Log entry on row 9 (D9 = 2)
# and appends a brand new Development Log entry on row 10:
#   Date = 3/4/2024, Time = 12:30 PM, Duration = 2, Task = "OO Design" (x2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Duration (h) value for the existing row 9 entry.
$ws.Range("D9").Value = 2

# Duplicate row 9's formatting down into row 10 so the new date/time cells
# reuse the workbook's existing date/time number-format styles instead of
# Excel auto-generating brand new ones.
$ws.Range("B9:F9").Copy()
$ws.Range("B10:F10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row 10 values.
$ws.Range("B10").Value = (Get-Date -Year 2024 -Month 3 -Day 4 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C10").Value = 0.52083333333333337
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "OO Design"
$ws.Range("F10").Value = "OO Design"

# Match the author's final selection (active cell moved to D10).
$ws.Range("D10").Select()
